$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '98.301.46'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = '3.419.16'
$ws.Range('E3').Value = '  +2.32%  '
$ws.Range('E4').Value = '  -0.01%  '
Set-TextValue 'D5' '255.84'
$ws.Range('E5').Value = '  -1.15%  '
Set-TextValue 'D6' '668.99'
$ws.Range('E6').Value = '  +1.67%  '
$ws.Range('E7').Value = '  -5.49%  '
Set-TextValue 'D8' '0.437'
$ws.Range('E8').Value = '  -4.56%  '
$ws.Range('E9').Value = '  -2.31%  '
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('D11').Value = '3.415.76'
$ws.Range('E11').Value = '  +2.35%  '
$ws.Range('E12').Value = '  +3.07%  '
Set-TextValue 'D13' '42.16'
$ws.Range('E13').Value = '  -2.19%  '
Set-TextValue 'D14' '6.43'
$ws.Range('E14').Value = '  +14.81%  '
$ws.Range('D15').Value = '98.083.97'
$ws.Range('E15').Value = '  -1.15%  '
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').Value = '4.049.96'
$ws.Range('E17').Value = '  +1.80%  '
Set-TextValue 'D18' '9.00'
$ws.Range('E18').Value = '  +19.71%  '
$ws.Range('D19').Value = '3.416.02'
$ws.Range('E19').Value = '  +2.08%  '
Set-TextValue 'D20' '0.580'
$ws.Range('E20').Value = '  +33.46%  '
Set-TextValue 'D21' '17.68'
$ws.Range('E21').Value = '  +4.76%  '
Set-TextValue 'D22' '11.05'
$ws.Range('E22').Value = '  +5.26%  '
Set-TextValue 'D23' '3.46'
$ws.Range('E23').Value = '  -4.60%  '
Set-TextValue 'D24' '512.06'
$ws.Range('E24').Value = '  -3.73%  '
Set-TextValue 'D25' '0.0000207'
Set-TextValue 'D26' '6.65'
$ws.Range('E26').Value = '  +6.09%  '
Set-TextValue 'D27' '101.78'
$ws.Range('E27').Value = '  +0.23%  '
Set-TextValue 'D28' '12.87'
$ws.Range('E28').Value = '  +2.21%  '
$ws.Range('D29').Value = '3.603.73'
$ws.Range('E29').Value = '  +2.17%  '
$ws.Range('E30').Value = '  +1.03%  '
Set-TextValue 'D31' '11.62'
$ws.Range('E31').Value = '  +5.53%  '
Set-TextValue 'D32' '0.198'
$ws.Range('E32').Value = '  +2.75%  '
Set-TextValue 'D33' '1.00'
$ws.Range('E33').Value = '  -0.02%  '
Set-TextValue 'D34' '2.53'
$ws.Range('E34').Value = '  +19.99%  '
Set-TextValue 'D35' '0.579'
$ws.Range('E35').Value = '  +7.75%  '
Set-TextValue 'D36' '0.998'
$ws.Range('E36').Value = '  -0.24%  '
Set-TextValue 'D37' '30.10'
$ws.Range('E37').Value = '  +2.57%  '
Set-TextValue 'D38' '1.52'
$ws.Range('E38').Value = '  +14.59%  '
Set-TextValue 'D39' '7.95'
$ws.Range('E39').Value = '  +1.25%  '
Set-TextValue 'D40' '538.88'
$ws.Range('E40').Value = '  +1.92%  '
$ws.Range('E41').Value = '  -2.73%  '
$ws.Range('E42').Value = '  -0.01%  '
Set-TextValue 'D43' '0.879'
$ws.Range('E43').Value = '  +6.46%  '
Set-TextValue 'D44' '24.72'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D45' '5.90'
$ws.Range('E45').Value = '  +15.17%  '
$ws.Range('B46').Value = 'Cosmos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D46' '9.05'
$ws.Range('E46').Value = '  +13.85%  '
$ws.Range('E47').Value = '  +0.58%  '
Set-TextValue 'D48' '3.78'
$ws.Range('E48').Value = '  +0.37%  '
Set-TextValue 'D49' '1.75'
$ws.Range('E49').Value = '  +16.58%  '
Set-TextValue 'D50' '3.28'
$ws.Range('E50').Value = '  -2.67%  '
Set-TextValue 'D51' '54.12'
$ws.Range('E51').Value = '  +10.13%  '
